$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New country/count data (updated search query results)
$data = @(
    ,@("PEOPLES R CHINA", 255)
    ,@("USA", 98)
    ,@("JAPAN", 63)
    ,@("FRANCE", 57)
    ,@("GERMANY", 44)
    ,@("NETHERLANDS", 39)
    ,@("INDIA", 31)
    ,@("ENGLAND", 28)
    ,@("RUSSIA", 19)
    ,@("SOUTH KOREA", 18)
    ,@("ITALY", 12)
    ,@("AUSTRALIA", 9)
    ,@("SPAIN", 9)
    ,@("TAIWAN", 8)
    ,@("IRAN", 7)
    ,@("PORTUGAL", 7)
    ,@("SAUDI ARABIA", 7)
    ,@("EGYPT", 6)
    ,@("CANADA", 5)
    ,@("PAKISTAN", 5)
    ,@("POLAND", 5)
    ,@("TURKEY", 5)
    ,@("BELGIUM", 4)
    ,@("BRAZIL", 4)
    ,@("SWITZERLAND", 4)
    ,@("UZBEKISTAN", 4)
    ,@("BANGLADESH", 3)
    ,@("IRAQ", 3)
    ,@("MALAYSIA", 3)
    ,@("MEXICO", 3)
    ,@("ROMANIA", 3)
    ,@("SWEDEN", 3)
    ,@("AUSTRIA", 2)
    ,@("CHILE", 2)
    ,@("DENMARK", 2)
    ,@("ISRAEL", 2)
    ,@("MOROCCO", 2)
    ,@("SCOTLAND", 2)
    ,@("SINGAPORE", 2)
    ,@("SLOVAKIA", 2)
    ,@("SOUTH AFRICA", 2)
    ,@("UKRAINE", 2)
    ,@("BELARUS", 1)
    ,@("BULGARIA", 1)
    ,@("COSTA RICA", 1)
    ,@("COTE IVOIRE", 1)
    ,@("CYPRUS", 1)
    ,@("CZECH REPUBLIC", 1)
    ,@("ETHIOPIA", 1)
    ,@("GEORGIA", 1)
    ,@("GREECE", 1)
    ,@("INDONESIA", 1)
    ,@("IRELAND", 1)
    ,@("JORDAN", 1)
    ,@("KAZAKHSTAN", 1)
    ,@("KENYA", 1)
    ,@("NIGERIA", 1)
    ,@("QATAR", 1)
    ,@("VIETNAM", 1)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = [double]$data[$i][1]
}
